# Apply the vocabulary.xlsx update: fill in format-description rows (123-130)
# and append new rows (131-142) for additional file-format vocabulary terms.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 123: add "format" label to D123 ---
$ws.Range("D123").Value = "format"

# --- Row 124: PDF/A ---
$ws.Range("B124").Value = "PDF/A - Portable Document Format for archiving "
$ws.Range("E124").Value = "A DANS preferred text format with '.pdf' ending."
$ws.Range("G124").Value = "https://dans.knaw.nl/en/file-formats/text-documents/pdf-a/"

# --- Row 125: PDF ---
$ws.Range("B125").Value = "PDF - Portable Document Format "
$ws.Range("E125").Value = "A DANS non-preferred text format, which includes all PDF other than PDF/A wiht '.pdf' ending."
$ws.Range("G125").Value = "https://dans.knaw.nl/en/file-formats/text-documents/pdf-pdf/"

# --- Row 126: ODT ---
$ws.Range("B126").Value = "ODT - Open Document "
$ws.Range("E126").Value = "A DANS preferred text format with '.odt' ending."
$ws.Range("G126").Value = "https://dans.knaw.nl/en/file-formats/text-documents/opendocument-text/"

# --- Row 127: DOC ---
$ws.Range("B127").Value = "DOC - Microsoft Word"
$ws.Range("E127").Value = "A DANS non-preferred text format with '.doc' ending."
$ws.Range("G127").Value = "https://dans.knaw.nl/en/file-formats/text-documents/microsoft-word-and-office-open-xml/"

# --- Row 128: DOCX ---
$ws.Range("B128").Value = "DOCX - Office Open XML"
$ws.Range("E128").Value = "A DANS non-preferred text format with '.docx' ending."
$ws.Range("G128").Value = "https://dans.knaw.nl/en/file-formats/text-documents/microsoft-word-and-office-open-xml/"

# --- Row 129: RTF ---
$ws.Range("B129").Value = "RTF - Rich Text File"
$ws.Range("E129").Value = "A DANS non-preferred text format with '.rtf' ending."
$ws.Range("G129").Value = "https://dans.knaw.nl/en/file-formats/text-documents/rich-text-file-rtf/"

# --- Row 130: TXT (Unicode) ---
$ws.Range("B130").Value = "TXT - Unicode text"
$ws.Range("E130").Value = "A DANS preferred plain text format with '.txt' ending."
$ws.Range("G130").Value = "https://dans.knaw.nl/en/file-formats/plain-text/unicode/"

# --- New rows 131-142 ---
# Columns: A=Identifier, B=skos:prefLabel@en, E=skos:definition@en, G=dct:source, H=skos:broader
$newRows = @(
    @{ Row = 131; A = "gen:10113"; B = "TXT - Non-Unicode text "; E = "A DANS non-preferred plain text format with '.txt' ending."; G = "https://dans.knaw.nl/en/file-formats/plain-text/unicode/"; H = "gen:10105" },
    @{ Row = 132; A = "gen:10114"; B = "XML - eXtensivle Markup Language"; E = "A DANS preferred markup language format with '.xml' ending."; G = "https://dans.knaw.nl/en/file-formats/markup-language/xml/"; H = "gen:10105" },
    @{ Row = 133; A = "gen:10115"; B = "HTML - Hypertext Markup Language"; E = "A DANS preferred markup language format with '.xml' ending."; G = "https://dans.knaw.nl/en/file-formats/markup-language/html/"; H = "gen:10105" },
    @{ Row = 134; A = "gen:10116"; B = "CSS - Cascading Style Sheets"; E = "A DANS preferred markup language format with '.css' ending."; G = "https://dans.knaw.nl/en/file-formats/markup-language/css/"; H = "gen:10105" },
    @{ Row = 135; A = "gen:10117"; B = "XSLT - Extensible Stylesheet Language Transformations"; E = "A DANS preferred markup language format with '.xslt' ending."; G = "https://dans.knaw.nl/en/file-formats/markup-language/xslt/"; H = "gen:10105" },
    @{ Row = 136; A = "gen:10118"; B = "JS - JavaScript"; E = "A DANS preferred markup language format with '.js' ending."; G = "https://dans.knaw.nl/en/file-formats/markup-language/script/"; H = "gen:10105" },
    @{ Row = 137; A = "gen:10119"; B = "ES - ECMAScript"; E = "A DANS preferred markup language format with '.es' ending."; G = "https://dans.knaw.nl/en/file-formats/markup-language/script/"; H = "gen:10105" },
    @{ Row = 138; A = "gen:10120"; B = "SGML - Standard Generalized Markup Language"; E = "A DANS non-preferred markup language format with '.es' ending."; G = "https://dans.knaw.nl/en/file-formats/markup-language/sgml/"; H = "gen:10105" },
    @{ Row = 139; A = "gen:10121"; B = "MD - Markdown"; E = "A DANS non-preferred markup language format with '.md' ending."; G = "https://dans.knaw.nl/en/file-formats/markup-language/markdown/"; H = "gen:10105" },
    @{ Row = 140; A = "gen:10122"; B = "MATLAB"; E = "A DANS preferred programming language format. "; G = "https://dans.knaw.nl/en/file-formats/programming-languages/matlab/"; H = "gen:10105" },
    @{ Row = 141; A = "gen:10123"; B = "NetCDF - Network Common Data Form"; E = "A DANS preferred file format. It is an interface to a library of data access functions for storing and retrieving data in the form of arrays."; G = "https://dans.knaw.nl/en/file-formats/programming-languages/netcdf/"; H = "gen:10105" },
    @{ Row = 142; A = "gen:10124"; B = "TF - Text-Fabric"; E = "A DANS preferred file format. Files with this format store a column of feature values that correspond to nodes and edges in a graph, which together represent annotated text. "; G = "https://dans.knaw.nl/en/file-formats/programming-languages/text-fabric/"; H = "gen:10105" }
)

foreach ($r in $newRows) {
    $ws.Range("A" + $r.Row).Value = $r.A
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("E" + $r.Row).Value = $r.E
    $ws.Range("G" + $r.Row).Value = $r.G
    $ws.Range("H" + $r.Row).Value = $r.H
}
